# "Add files via upload" -- adds three new "*_arrow" columns (D:F) to the
# "axes" sheet, mirroring the existing A/B/C (FeO / Na2O+K2O / MgO) values,
# and shifts the old "Title/AFM" column out to G. Also leaves the "axes"
# sheet as the active/selected tab (it was the last one touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("axes")

# --- shift the old D column (header + value) out to the new G column ---
$ws.Range("G1").Value = $ws.Range("D1").Value()
$ws.Range("G2").Value = $ws.Range("D2").Value()

# --- new header row (D1:F1) ---
$ws.Range("D1").Value = "A_arrow"
$ws.Range("E1").Value = "B_arrow"
$ws.Range("F1").Value = "C_arrow"

# --- new data row (D2:F2) mirrors A2:C2 (FeO / Na2O + K2O / MgO) ---
$ws.Range("D2").Value = $ws.Range("A2").Value()
$ws.Range("E2").Value = $ws.Range("B2").Value()
$ws.Range("F2").Value = $ws.Range("C2").Value()

# widen the three new columns
$ws.Columns("D:F").ColumnWidth = 13.25

# make "axes" the active sheet / tab, with the new selection
$ws.Activate()
$ws.Range("B5").Select()
